$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2026-02-02 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-02-03 Tuesday", 2) | Out-Null

# Update the division-problem table, cell by cell (row, col), since
# some old values repeat with different replacements depending on position.
$t = $d.Tables(1)

$t.Cell(1, 1).Range.Text = "25÷5=5, 0"
$t.Cell(1, 2).Range.Text = "72÷6=12, 0"
$t.Cell(1, 3).Range.Text = "38÷8=4, 6"
$t.Cell(1, 4).Range.Text = "98÷6=16, 2"
$t.Cell(1, 5).Range.Text = "61÷4=15, 1"
$t.Cell(5, 1).Range.Text = "74÷4=18, 2"
$t.Cell(5, 2).Range.Text = "73÷9=8, 1"
$t.Cell(5, 3).Range.Text = "28÷2=14, 0"
$t.Cell(5, 4).Range.Text = "62÷6=10, 2"
$t.Cell(5, 5).Range.Text = "35÷6=5, 5"
$t.Cell(9, 1).Range.Text = "60÷6=10, 0"
$t.Cell(9, 2).Range.Text = "79÷4=19, 3"
$t.Cell(9, 3).Range.Text = "54÷2=27, 0"
$t.Cell(9, 4).Range.Text = "71÷4=17, 3"
$t.Cell(9, 5).Range.Text = "80÷9=8, 8"
$t.Cell(13, 1).Range.Text = "81÷9=9, 0"
$t.Cell(13, 2).Range.Text = "85÷2=42, 1"
$t.Cell(13, 3).Range.Text = "50÷9=5, 5"
$t.Cell(13, 4).Range.Text = "24÷3=8, 0"
$t.Cell(13, 5).Range.Text = "32÷6=5, 2"
$t.Cell(17, 1).Range.Text = "85÷9=9, 4"
$t.Cell(17, 2).Range.Text = "61÷3=20, 1"
$t.Cell(17, 3).Range.Text = "20÷5=4, 0"
$t.Cell(17, 4).Range.Text = "84÷6=14, 0"
$t.Cell(17, 5).Range.Text = "83÷5=16, 3"
